$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 loses its A:H data (word no longer present in this slot)
$ws.Range("A11:H11").Clear()

# Row 1
$ws.Range("A1").Value = 'negative'
$ws.Range("J1").Value = 'positive'

# Row 2
$ws.Range("A2").Value = 'name'
$ws.Range("B2").Value = 'anchor score'
$ws.Range("C2").Value = 'type occurences'
$ws.Range("D2").Value = 'total occurences'
$ws.Range("E2").Value = '+%'
$ws.Range("F2").Value = '-%'
$ws.Range("G2").Value = 'both'
$ws.Range("H2").Value = 'normal'
$ws.Range("J2").Value = 'name'
$ws.Range("K2").Value = 'anchor score'
$ws.Range("L2").Value = 'type occurences'
$ws.Range("M2").Value = 'total occurences'
$ws.Range("N2").Value = '+%'
$ws.Range("O2").Value = '-%'
$ws.Range("P2").Value = 'both'
$ws.Range("Q2").Value = 'normal'

# Row 3
$ws.Range("A3").Value = 'died'
$ws.Range("B3").Value = 0.7307692307692307
$ws.Range("C3").Value = 19
$ws.Range("D3").Value = 19
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 7
$ws.Range("J3").Value = 'interesting'
$ws.Range("K3").Value = 0.9393939393939394
$ws.Range("L3").Value = 31
$ws.Range("M3").Value = 31
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 2

# Row 4
$ws.Range("A4").Value = 'crude'
$ws.Range("B4").Value = 0.7058823529411765
$ws.Range("C4").Value = 24
$ws.Range("D4").Value = 24
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 10
$ws.Range("J4").Value = 'best'
$ws.Range("K4").Value = 0.9152542372881356
$ws.Range("L4").Value = 54
$ws.Range("M4").Value = 54
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 5

# Row 5
$ws.Range("A5").Value = 'forced'
$ws.Range("B5").Value = 0.6551724137931034
$ws.Range("C5").Value = 19
$ws.Range("D5").Value = 19
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 10
$ws.Range("J5").Value = 'love'
$ws.Range("K5").Value = 0.8913043478260869
$ws.Range("L5").Value = 41
$ws.Range("M5").Value = 41
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 5

# Row 6
$ws.Range("A6").Value = 'fraud'
$ws.Range("B6").Value = 0.6111111111111112
$ws.Range("C6").Value = 22
$ws.Range("D6").Value = 22
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 14
$ws.Range("J6").Value = 'nice'
$ws.Range("K6").Value = 0.8888888888888888
$ws.Range("L6").Value = 24
$ws.Range("M6").Value = 24
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 3

# Row 7
$ws.Range("A7").Value = 'crisis'
$ws.Range("B7").Value = 0.5958904109589042
$ws.Range("C7").Value = 174
$ws.Range("D7").Value = 174
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 118
$ws.Range("J7").Value = 'great'
$ws.Range("K7").Value = 0.8482142857142857
$ws.Range("L7").Value = 95
$ws.Range("M7").Value = 95
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 17

# Row 8
$ws.Range("A8").Value = 'panic'
$ws.Range("B8").Value = 0.2093023255813954
$ws.Range("C8").Value = 108
$ws.Range("D8").Value = 108
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 408
$ws.Range("J8").Value = 'happy'
$ws.Range("K8").Value = 0.8461538461538461
$ws.Range("L8").Value = 22
$ws.Range("M8").Value = 22
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 4

# Row 9
$ws.Range("A9").Value = 'sc'
$ws.Range("B9").Value = 0.2063492063492063
$ws.Range("C9").Value = 39
$ws.Range("D9").Value = 39
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = 150
$ws.Range("J9").Value = 'positive'
$ws.Range("K9").Value = 0.7931034482758621
$ws.Range("L9").Value = 46
$ws.Range("M9").Value = 46
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 12

# Row 10
$ws.Range("A10").Value = 'low'
$ws.Range("B10").Value = 0.1409395973154362
$ws.Range("C10").Value = 21
$ws.Range("D10").Value = 21
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = $false
$ws.Range("H10").Value = 128
$ws.Range("J10").Value = 'thank'
$ws.Range("K10").Value = 0.78125
$ws.Range("L10").Value = 100
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 28

# Row 11
$ws.Range("J11").Value = 'special'
$ws.Range("K11").Value = 0.7777777777777778
$ws.Range("L11").Value = 28
$ws.Range("M11").Value = 28
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 8

# Row 12
$ws.Range("J12").Value = 'healthy'
$ws.Range("K12").Value = 0.7777777777777778
$ws.Range("L12").Value = 21
$ws.Range("M12").Value = 21
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 6

# Row 13
$ws.Range("J13").Value = 'thanks'
$ws.Range("K13").Value = 0.7682926829268293
$ws.Range("L13").Value = 63
$ws.Range("M13").Value = 63
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 19

# Row 14
$ws.Range("J14").Value = 'free'
$ws.Range("K14").Value = 0.7333333333333333
$ws.Range("L14").Value = 88
$ws.Range("M14").Value = 88
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 32

# Row 15
$ws.Range("J15").Value = 'safe'
$ws.Range("K15").Value = 0.7323943661971831
$ws.Range("L15").Value = 104
$ws.Range("M15").Value = 104
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 38

# Row 16
$ws.Range("J16").Value = 'safety'
$ws.Range("K16").Value = 0.7058823529411765
$ws.Range("L16").Value = 36
$ws.Range("M16").Value = 36
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 15

# Row 17
$ws.Range("J17").Value = 'good'
$ws.Range("K17").Value = 0.7
$ws.Range("L17").Value = 112
$ws.Range("M17").Value = 112
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 48

# Row 18
$ws.Range("J18").Value = 'support'
$ws.Range("K18").Value = 0.6886792452830188
$ws.Range("L18").Value = 73
$ws.Range("M18").Value = 73
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 33

# Row 19
$ws.Range("J19").Value = 'heroes'
$ws.Range("K19").Value = 0.6595744680851063
$ws.Range("L19").Value = 31
$ws.Range("M19").Value = 31
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 16

# Row 20
$ws.Range("J20").Value = 'confidence'
$ws.Range("K20").Value = 0.6388888888888888
$ws.Range("L20").Value = 23
$ws.Range("M20").Value = 23
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 13

# Row 21
$ws.Range("J21").Value = 'better'
$ws.Range("K21").Value = 0.6031746031746031
$ws.Range("L21").Value = 38
$ws.Range("M21").Value = 38
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = $false
$ws.Range("Q21").Value = 25

# Row 22
$ws.Range("J22").Value = 'well'
$ws.Range("K22").Value = 0.5851063829787234
$ws.Range("L22").Value = 55
$ws.Range("M22").Value = 55
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = $false
$ws.Range("Q22").Value = 39

# Row 23
$ws.Range("J23").Value = 'fresh'
$ws.Range("K23").Value = 0.5833333333333334
$ws.Range("L23").Value = 28
$ws.Range("M23").Value = 28
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = $false
$ws.Range("Q23").Value = 20

# Row 24
$ws.Range("J24").Value = 'relief'
$ws.Range("K24").Value = 0.58
$ws.Range("L24").Value = 29
$ws.Range("M24").Value = 29
$ws.Range("N24").Value = 1
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = $false
$ws.Range("Q24").Value = 21

# Row 25
$ws.Range("J25").Value = 'credit'
$ws.Range("K25").Value = 0.5588235294117647
$ws.Range("L25").Value = 19
$ws.Range("M25").Value = 19
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = $false
$ws.Range("Q25").Value = 15

# Row 26
$ws.Range("J26").Value = 'important'
$ws.Range("K26").Value = 0.5333333333333333
$ws.Range("L26").Value = 24
$ws.Range("M26").Value = 24
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = $false
$ws.Range("Q26").Value = 21

# Row 27
$ws.Range("J27").Value = 'hand'
$ws.Range("K27").Value = 0.5195822454308094
$ws.Range("L27").Value = 199
$ws.Range("M27").Value = 199
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = $false
$ws.Range("Q27").Value = 184

# Row 28
$ws.Range("J28").Value = 'like'
$ws.Range("K28").Value = 0.4735294117647059
$ws.Range("L28").Value = 161
$ws.Range("M28").Value = 161
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = $false
$ws.Range("Q28").Value = 179

# Row 29
$ws.Range("J29").Value = 'help'
$ws.Range("K29").Value = 0.464406779661017
$ws.Range("L29").Value = 137
$ws.Range("M29").Value = 137
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = $false
$ws.Range("Q29").Value = 158

# Row 30
$ws.Range("J30").Value = 'care'
$ws.Range("K30").Value = 0.4382022471910113
$ws.Range("L30").Value = 39
$ws.Range("M30").Value = 39
$ws.Range("N30").Value = 1
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = $false
$ws.Range("Q30").Value = 50

# Row 31
$ws.Range("J31").Value = 'hope'
$ws.Range("K31").Value = 0.4
$ws.Range("L31").Value = 26
$ws.Range("M31").Value = 26
$ws.Range("N31").Value = 1
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = $false
$ws.Range("Q31").Value = 39

# Row 32
$ws.Range("J32").Value = 'increase'
$ws.Range("K32").Value = 0.3974358974358974
$ws.Range("L32").Value = 31
$ws.Range("M32").Value = 31
$ws.Range("N32").Value = 1
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = $false
$ws.Range("Q32").Value = 47

# Row 33
$ws.Range("J33").Value = 'protect'
$ws.Range("K33").Value = 0.3972602739726027
$ws.Range("L33").Value = 29
$ws.Range("M33").Value = 29
$ws.Range("N33").Value = 1
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = $false
$ws.Range("Q33").Value = 44

# Row 34
$ws.Range("J34").Value = 'sure'
$ws.Range("K34").Value = 0.34375
$ws.Range("L34").Value = 22
$ws.Range("M34").Value = 22
$ws.Range("N34").Value = 1
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = $false
$ws.Range("Q34").Value = 42

# Row 35
$ws.Range("J35").Value = 'please'
$ws.Range("K35").Value = 0.3430962343096234
$ws.Range("L35").Value = 82
$ws.Range("M35").Value = 82
$ws.Range("N35").Value = 1
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = $false
$ws.Range("Q35").Value = 157

# Row 36
$ws.Range("J36").Value = 'share'
$ws.Range("K36").Value = 0.2857142857142857
$ws.Range("L36").Value = 20
$ws.Range("M36").Value = 20
$ws.Range("N36").Value = 1
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = $false
$ws.Range("Q36").Value = 50

# Row 37
$ws.Range("J37").Value = 'store'
$ws.Range("K37").Value = 0.04026845637583892
$ws.Range("L37").Value = 36
$ws.Range("M37").Value = 36
$ws.Range("N37").Value = 1
$ws.Range("O37").Value = 0
$ws.Range("P37").Value = $false
$ws.Range("Q37").Value = 858

# Row 38
$ws.Range("J38").Value = 'grocery'
$ws.Range("K38").Value = 0.03329633740288569
$ws.Range("L38").Value = 30
$ws.Range("M38").Value = 30
$ws.Range("N38").Value = 1
$ws.Range("O38").Value = 0
$ws.Range("P38").Value = $false
$ws.Range("Q38").Value = 871

# Row 39
$ws.Range("J39").Value = '19'
$ws.Range("K39").Value = 0.01261682242990654
$ws.Range("L39").Value = 27
$ws.Range("M39").Value = 30
$ws.Range("N39").Value = 0.9
$ws.Range("O39").Value = 0.09999999999999998
$ws.Range("P39").Value = $true
$ws.Range("Q39").Value = 2113

# Row 40
$ws.Range("J40").Value = 'co'
$ws.Range("K40").Value = 0.01095713825330326
$ws.Range("L40").Value = 34
$ws.Range("M40").Value = 38
$ws.Range("N40").Value = 0.89
$ws.Range("O40").Value = 0.11
$ws.Range("P40").Value = $true
$ws.Range("Q40").Value = 3069

# Row 41
$ws.Range("J41").Value = 'corona'
$ws.Range("K41").Value = 0.006259780907668232
$ws.Range("L41").Value = 20
$ws.Range("M41").Value = 26
$ws.Range("N41").Value = 0.77
$ws.Range("O41").Value = 0.23
$ws.Range("P41").Value = $true
$ws.Range("Q41").Value = 3175
